# Insert a new data row at row 227 (pushes existing rows 227-274 down to 228-275)
# and populate it with the new record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 227, shifting rows 227..274 down to 228..275.
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new record's data.
$row = 227

$ws.Cells.Item($row, 1).Value2  = 5
$ws.Cells.Item($row, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value2  = "Maule"
$ws.Cells.Item($row, 4).Value2  = 45204
$ws.Cells.Item($row, 5).Value2  = 7
$ws.Cells.Item($row, 6).Value2  = 100112031
$ws.Cells.Item($row, 7).Value2  = "Poroto verde"
$ws.Cells.Item($row, 8).Value2  = "Sin especificar"
$ws.Cells.Item($row, 9).Value2  = "Primera"
$ws.Cells.Item($row, 10).Value2 = 150
$ws.Cells.Item($row, 11).Value2 = 23000
$ws.Cells.Item($row, 12).Value2 = 23000
$ws.Cells.Item($row, 13).Value2 = 23000
$ws.Cells.Item($row, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item($row, 15).Value2 = "Perú"
$ws.Cells.Item($row, 16).Value2 = 920
$ws.Cells.Item($row, 17).Value2 = 25
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
